$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Feature rows: row number, feature name, Timepoint1, Timepoint2, Percentage Change
# "Total_Brain_Volume" is inserted as a new feature (rows 2 and 4), and the
# existing "Lesions" row (previously the last row) moves up to row 3; all
# other feature rows shift down by 3 rows to make room.
$data = @(
    @(2, "Total_Brain_Volume", 702055.389495, 692432.848044, -1.4),
    @(3, "Lesions", 20861.282386, 20612.32408, -1.2),
    @(4, "Total_Brain_Volume", 702055.389495, 692432.848044, -1.4),
    @(5, "Left_Cerebral_Cortex", 216363.901702, 210382.948221, -2.8),
    @(6, "Right_Cerebral_Cortex", 217479.315374, 211679.885892, -2.7),
    @(7, "Left_Cerebral_White_Matter", 129820.814226, 131305.87288, 1.1),
    @(8, "Right_Cerebral_White_Matter", 138391.358193, 139064.141051, 0.5),
    @(9, "Left_Cerebellum_Cortex", 43376.990916, 42509.055112, -2),
    @(10, "Right_Cerebellum_Cortex", 42489.206922, 41462.946682, -2.4),
    @(11, "Left_Cerebellum_White_Matter", 12952.079447, 13143.855066, 1.5),
    @(12, "Right_Cerebellum_White_Matter", 12188.86018, 12389.316821, 1.6),
    @(13, "Left_Hippocampus", 3384.294242, 3320.398542, -1.9),
    @(14, "Right_Hippocampus", 3973.598266, 3876.511076, -2.4),
    @(15, "Left_Amygdala", 1541.30063, 1532.384809, -0.6),
    @(16, "Right_Amygdala", 1672.987502, 1637.421468, -2.1),
    @(17, "Left_VentralDC", 2882.265225, 2877.466901, -0.2),
    @(18, "Right_VentralDC", 2708.488099, 2696.681429, -0.4),
    @(19, "Left_Putamen", 4078.4353, 4087.333629, 0.2),
    @(20, "Right_Putamen", 3950.083621, 3961.549178, 0.3),
    @(21, "Left_Accumbens_area", 552.470008, 537.3571020000001, -2.7),
    @(22, "Right_Accumbens_area", 543.452117, 532.985835, -1.9),
    @(23, "Brain_Stem", 14719.269037, 14541.180806, -1.2),
    @(24, "Right_Pallidum", 2307.10325, 2304.277791, -0.1),
    @(25, "Left_Caudate", 2807.284109, 2705.534508, -3.6),
    @(26, "Right_Thalamus", 5078.785924, 5080.294255, 0),
    @(27, "Left_Pallidum", 1904.840295, 1929.677608, 1.3),
    @(28, "Right_Caudate", 2762.606931, 2681.386482, -2.9),
    @(29, "Left_Thalamus", 4548.971697, 4540.235741, -0.2),
    @(30, "Left_Lateral_Ventricle", 29540.644166, 30892.433138, 4.6),
    @(31, "Right_Lateral_Ventricle", 26514.731298, 27741.226193, 4.6),
    @(32, "Left_Inf_Lat_Vent", 858.196462, 940.390008, 9.6),
    @(33, "Right_Inf_Lat_Vent", 694.158335, 774.1007540000001, 11.5),
    @(34, "x3rd_Ventricle", 1993.551358, 2052.609872, 3),
    @(35, "x4th_Ventricle", 1192.331951, 1172.213845, -1.7),
    @(36, "x5th_Ventricle", 6.896968, 6.687121, -3),
    @(37, "CSF", 389216.034081, 400722.91654, 3)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
}

# Rows 36 and 37 are brand-new rows beyond the original sheet's used range
# (A1:D35), so column A needs the same bold/centered/bordered "feature name"
# formatting used by every other row in column A. Copy it from row 34 (an
# existing, untouched feature row) onto the new rows.
$ws.Range("A34").Copy() | Out-Null
$ws.Range("A36:A37").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
